$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 413785.2
$ws.Range("J17").Value = 441346.34
$ws.Range("L17").Value = 1324039.02
$ws.Range("N17").Value = -1324375.02
$ws.Range("H80").Value = 2329.375
$ws.Range("I80").Value = 4326.5
$ws.Range("J80").Value = 332.25
$ws.Range("K80").Value = 12979.5
$ws.Range("L80").Value = 996.75
$ws.Range("M80").Value = -11981.5
$ws.Range("N80").Value = -2992.75
$ws.Range("H83").Value = 2329.375
$ws.Range("I83").Value = 4326.5
$ws.Range("J83").Value = 332.25
$ws.Range("K83").Value = 38938.5
$ws.Range("L83").Value = 2990.25
$ws.Range("M83").Value = -33946.5
$ws.Range("N83").Value = -12974.25
$ws.Range("H92").Value = 798.1
$ws.Range("I92").Value = 814
$ws.Range("J92").Value = 655
$ws.Range("K92").Value = 814
$ws.Range("L92").Value = 655
$ws.Range("M92").Value = 434
$ws.Range("N92").Value = -3151
$ws.Range("H113").Value = 3822.3635
$ws.Range("I113").Value = 4099
$ws.Range("J113").Value = 3084.6667
$ws.Range("K113").Value = 4099
$ws.Range("L113").Value = 3084.6667
$ws.Range("M113").Value = -845
$ws.Range("N113").Value = -9592.6667
$ws.Range("H118").Value = 2089.8333
$ws.Range("I118").Value = 2089.8333
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 6269.499899999999
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -4612.499899999999
$ws.Range("N118").Value = ""
$ws.Range("H132").Value = 1831.4667
$ws.Range("I132").Value = 1864.2142
$ws.Range("J132").Value = 1373
$ws.Range("K132").Value = 5592.642599999999
$ws.Range("L132").Value = 4119
$ws.Range("M132").Value = -3062.642599999999
$ws.Range("N132").Value = -9179
$ws.Range("H134").Value = 107250
$ws.Range("J134").Value = 107250
$ws.Range("L134").Value = 107250
$ws.Range("N134").Value = -117390
$ws.Range("H137").Value = 2313.7917
$ws.Range("I137").Value = 1534.3529
$ws.Range("K137").Value = 4603.0587
$ws.Range("M137").Value = -2053.0587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38464640
$ws.Range("I61").Value = 62501630
$ws.Range("K61").Value = 62501630
$ws.Range("M61").Value = -62501418
$ws.Range("H74").Value = 62502504
$ws.Range("I74").Value = 62502504
$ws.Range("K74").Value = 62502504
$ws.Range("M74").Value = -62501630
$ws.Range("H77").Value = 62502504
$ws.Range("I77").Value = 62502504
$ws.Range("K77").Value = 312512520
$ws.Range("M77").Value = -312508152
$ws.Range("H102").Value = 1732.65
$ws.Range("I102").Value = 1270.5834
$ws.Range("J102").Value = 2425.75
$ws.Range("K102").Value = 1270.5834
$ws.Range("L102").Value = 2425.75
$ws.Range("M102").Value = 351.4166
$ws.Range("N102").Value = -5669.75
$ws.Range("H132").Value = 2704754.2
$ws.Range("I132").Value = 3227513.5
$ws.Range("J132").Value = 3831.1667
$ws.Range("K132").Value = 9682540.5
$ws.Range("L132").Value = 11493.5001
$ws.Range("M132").Value = -9680010.5
$ws.Range("N132").Value = -16553.5001
$ws.Range("H136").Value = 38464640
$ws.Range("I136").Value = 62501630
$ws.Range("K136").Value = 187504890
$ws.Range("M136").Value = -187502340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 16455300
$ws.Range("I134").Value = 17003676
$ws.Range("K134").Value = 51011028
$ws.Range("M134").Value = -51008493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2956.7778
$ws.Range("I7").Value = 3230.1428
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 3230.1428
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -3117.1428
$ws.Range("N7").Value = -2226
$ws.Range("H31").Value = 4100.222
$ws.Range("I31").Value = 2838.4285
$ws.Range("J31").Value = 8516.5
$ws.Range("K31").Value = 2838.4285
$ws.Range("L31").Value = 8516.5
$ws.Range("M31").Value = -2543.4285
$ws.Range("N31").Value = -9106.5
$ws.Range("H34").Value = 4100.222
$ws.Range("I34").Value = 2838.4285
$ws.Range("J34").Value = 8516.5
$ws.Range("K34").Value = 2838.4285
$ws.Range("L34").Value = 8516.5
$ws.Range("M34").Value = -2636.4285
$ws.Range("N34").Value = -8920.5
$ws.Range("H58").Value = 26323028
$ws.Range("I58").Value = 62513348
$ws.Range("K58").Value = 62513348
$ws.Range("M58").Value = -62513145
$ws.Range("H62").Value = 5610.857
$ws.Range("I62").Value = 5876.6
$ws.Range("J62").Value = 4946.5
$ws.Range("K62").Value = 5876.6
$ws.Range("L62").Value = 4946.5
$ws.Range("M62").Value = -5252.6
$ws.Range("N62").Value = -6194.5
$ws.Range("H65").Value = 5610.857
$ws.Range("I65").Value = 5876.6
$ws.Range("J65").Value = 4946.5
$ws.Range("K65").Value = 29383
$ws.Range("L65").Value = 24732.5
$ws.Range("M65").Value = -26263
$ws.Range("N65").Value = -30972.5
$ws.Range("H88").Value = 17199.334
$ws.Range("J88").Value = 18958.8
$ws.Range("L88").Value = 18958.8
$ws.Range("N88").Value = -19770.8
$ws.Range("H91").Value = 17199.334
$ws.Range("J91").Value = 18958.8
$ws.Range("L91").Value = 18958.8
$ws.Range("N91").Value = -21766.8
$ws.Range("H107").Value = 92009.55
$ws.Range("I107").Value = 502
$ws.Range("K107").Value = 502
$ws.Range("M107").Value = 1418
$ws.Range("H132").Value = 32260402
$ws.Range("I132").Value = 34484980
$ws.Range("K132").Value = 103454940
$ws.Range("M132").Value = -103452410
$ws.Range("H134").Value = 10042000
$ws.Range("I134").Value = 11954042
$ws.Range("J134").Value = 3777.25
$ws.Range("K134").Value = 35862126
$ws.Range("L134").Value = 11331.75
$ws.Range("M134").Value = -35859591
$ws.Range("N134").Value = -16401.75
$ws.Range("H136").Value = 26323028
$ws.Range("I136").Value = 62513348
$ws.Range("K136").Value = 187540044
$ws.Range("M136").Value = -187537494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 668701.7
$ws.Range("I4").Value = 858330.7
$ws.Range("K4").Value = 2574992.1
$ws.Range("M4").Value = -2574880.1
$ws.Range("H5").Value = 168833.17
$ws.Range("I5").Value = 501499.5
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 1504498.5
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -1504386.5
$ws.Range("N5").Value = -7724
$ws.Range("H29").Value = 5086
$ws.Range("I29").Value = 10000
$ws.Range("J29").Value = 172
$ws.Range("K29").Value = 30000
$ws.Range("L29").Value = 516
$ws.Range("M29").Value = -29723
$ws.Range("N29").Value = -1070
$ws.Range("H76").Value = 15000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14617
$ws.Range("H79").Value = 15000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13674
$ws.Range("H112").Value = 12731.889
$ws.Range("J112").Value = 14855.286
$ws.Range("L112").Value = 44565.858
$ws.Range("N112").Value = -46781.858
$ws.Range("H135").Value = 168833.17
$ws.Range("I135").Value = 501499.5
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 4513495.5
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -4510960.5
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 30038900
$ws.Range("I11").Value = 37545000
$ws.Range("J11").Value = 14500
$ws.Range("K11").Value = 37545000
$ws.Range("L11").Value = 14500
$ws.Range("M11").Value = -37544861
$ws.Range("N11").Value = -14778
$ws.Range("H132").Value = 20835136
$ws.Range("I132").Value = 20835136
$ws.Range("K132").Value = 62505408
$ws.Range("M132").Value = -62502878

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 986.6667
$ws.Range("I16").Value = 966.1818
$ws.Range("K16").Value = 966.1818
$ws.Range("M16").Value = -796.1818
$ws.Range("H46").Value = 1955.4286
$ws.Range("I46").Value = 1955.4286
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1955.4286
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1767.4286
$ws.Range("N46").Value = ""
$ws.Range("H55").Value = 649.13635
$ws.Range("J55").Value = 811.4545
$ws.Range("L55").Value = 811.4545
$ws.Range("N55").Value = -1157.4545
$ws.Range("H56").Value = 18833
$ws.Range("J56").Value = 18833
$ws.Range("L56").Value = 18833
$ws.Range("N56").Value = -20215
$ws.Range("H100").Value = 15913199
$ws.Range("I100").Value = 17504050
$ws.Range("J100").Value = 4697
$ws.Range("K100").Value = 17504050
$ws.Range("L100").Value = 4697
$ws.Range("M100").Value = -17503509
$ws.Range("N100").Value = -5779
$ws.Range("H122").Value = 5370.364
$ws.Range("I122").Value = 3321.75
$ws.Range("K122").Value = 9965.25
$ws.Range("M122").Value = -7515.25
$ws.Range("H136").Value = 2076.625
$ws.Range("I136").Value = 2076.625
$ws.Range("K136").Value = 6229.875
$ws.Range("M136").Value = -3679.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 799.5
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 799.5
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 799.5
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -1025.5
$ws.Range("H122").Value = 1804.5555
$ws.Range("I122").Value = 1804.5555
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5413.666499999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2963.666499999999
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 62519464
$ws.Range("I132").Value = 100000740
$ws.Range("K132").Value = 300002220
$ws.Range("M132").Value = -299999690
$ws.Range("H135").Value = 65350
$ws.Range("J135").Value = 65350
$ws.Range("L135").Value = 65350
$ws.Range("N135").Value = -75490
